$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5725316666666668
$ws.Range("H2").Value = 1.717595
$ws.Range("I2").Value = 0.3864899584549088
$ws.Range("J2").Value = 0.3864899584549088
$ws.Range("M2").Value = 46.29121633333333
$ws.Range("N2").Value = 138.873649
$ws.Range("O2").Value = 0.3133663986859022
$ws.Range("P2").Value = 0.3133663986859022
$ws.Range("Q2").Value = 26.50318723935056
$ws.Range("R2").Value = 238.528685154155
$ws.Range("S2").Value = 0.1211129664092787
$ws.Range("T2").Value = 0.1211129664092787

# Row 3
$ws.Range("G3").Value = 0.5725316666666668
$ws.Range("H3").Value = 1.717595
$ws.Range("I3").Value = 0.3864899584549088
$ws.Range("J3").Value = 0.3864899584549088
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("O3").Value = 0.3169204109998198
$ws.Range("P3").Value = 0.3169204109998198
$ws.Range("Q3").Value = 26.80377037207223
$ws.Range("R3").Value = 241.23393334865
$ws.Range("S3").Value = 0.122486556480833
$ws.Range("T3").Value = 0.122486556480833

# Row 4
$ws.Range("G4").Value = 0.5725316666666668
$ws.Range("H4").Value = 1.717595
$ws.Range("I4").Value = 0.3864899584549088
$ws.Range("J4").Value = 0.3864899584549088
$ws.Range("M4").Value = 38.53544233333333
$ws.Range("N4").Value = 115.606327
$ws.Range("O4").Value = 0.2608640200510233
$ws.Range("P4").Value = 0.2608640200510233
$ws.Range("Q4").Value = 22.06276102484056
$ws.Range("R4").Value = 198.564849223565
$ws.Range("S4").Value = 0.1008213242719005
$ws.Range("T4").Value = 0.1008213242719005

# Row 5
$ws.Range("G5").Value = 0.5725316666666668
$ws.Range("H5").Value = 1.717595
$ws.Range("I5").Value = 0.3864899584549088
$ws.Range("J5").Value = 0.3864899584549088
$ws.Range("M5").Value = 16.07945366666667
$ws.Range("N5").Value = 48.238361
$ws.Range("O5").Value = 0.1088491702632547
$ws.Range("P5").Value = 0.1088491702632547
$ws.Range("Q5").Value = 9.205996406866113
$ws.Range("R5").Value = 82.85396766179501
$ws.Range("S5").Value = 0.04206911129289662
$ws.Range("T5").Value = 0.04206911129289662

# Row 6
$ws.Range("H6").Value = 0.919331
$ws.Range("I6").Value = 0.2068661122070742
$ws.Range("J6").Value = 0.2068661122070743
$ws.Range("M6").Value = 46.29121633333333
$ws.Range("N6").Value = 138.873649
$ws.Range("O6").Value = 0.3133663986859022
$ws.Range("P6").Value = 0.3133663986859022
$ws.Range("Q6").Value = 14.18565006764656
$ws.Range("R6").Value = 127.670850608819
$ws.Range("S6").Value = 0.06482488859248461
$ws.Range("T6").Value = 0.06482488859248461

# Row 7
$ws.Range("H7").Value = 0.919331
$ws.Range("I7").Value = 0.2068661122070742
$ws.Range("J7").Value = 0.2068661122070743
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("O7").Value = 0.3169204109998198
$ws.Range("P7").Value = 0.3169204109998198
$ws.Range("Q7").Value = 14.34653513775222
$ws.Range("S7").Value = 0.0655600933026008
$ws.Range("T7").Value = 0.0655600933026008

# Row 8
$ws.Range("H8").Value = 0.919331
$ws.Range("I8").Value = 0.2068661122070742
$ws.Range("J8").Value = 0.2068661122070743
$ws.Range("M8").Value = 38.53544233333333
$ws.Range("N8").Value = 115.606327
$ws.Range("O8").Value = 0.2608640200510233
$ws.Range("P8").Value = 0.2608640200510233
$ws.Range("Q8").Value = 11.80894224524855
$ws.Range("R8").Value = 106.280480207237
$ws.Range("S8").Value = 0.05396392564266345
$ws.Range("T8").Value = 0.05396392564266346

# Row 9
$ws.Range("H9").Value = 0.919331
$ws.Range("I9").Value = 0.2068661122070742
$ws.Range("J9").Value = 0.2068661122070743
$ws.Range("M9").Value = 16.07945366666667
$ws.Range("N9").Value = 48.238361
$ws.Range("O9").Value = 0.1088491702632547
$ws.Range("P9").Value = 0.1088491702632547
$ws.Range("Q9").Value = 4.927446739610111
$ws.Range("R9").Value = 44.347020656491
$ws.Range("S9").Value = 0.02251720466932538
$ws.Range("T9").Value = 0.02251720466932539

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5451493333333333
$ws.Range("H10").Value = 1.635448
$ws.Range("I10").Value = 0.3680053968340403
$ws.Range("J10").Value = 0.3680053968340404
$ws.Range("M10").Value = 46.29121633333333
$ws.Range("N10").Value = 138.873649
$ws.Range("O10").Value = 0.3133663986859022
$ws.Range("P10").Value = 0.3133663986859022
$ws.Range("Q10").Value = 25.23562572330577
$ws.Range("R10").Value = 227.120631509752
$ws.Range("S10").Value = 0.1153205259028595
$ws.Range("T10").Value = 0.1153205259028596

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5451493333333333
$ws.Range("H11").Value = 1.635448
$ws.Range("I11").Value = 0.3680053968340403
$ws.Range("J11").Value = 0.3680053968340404
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("O11").Value = 0.3169204109998198
$ws.Range("P11").Value = 0.3169204109998198
$ws.Range("Q11").Value = 25.52183293935111
$ws.Range("R11").Value = 229.69649645416
$ws.Range("S11").Value = 0.1166284216147958
$ws.Range("T11").Value = 0.1166284216147958

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5451493333333333
$ws.Range("H12").Value = 1.635448
$ws.Range("I12").Value = 0.3680053968340403
$ws.Range("J12").Value = 0.3680053968340404
$ws.Range("M12").Value = 38.53544233333333
$ws.Range("N12").Value = 115.606327
$ws.Range("O12").Value = 0.2608640200510233
$ws.Range("P12").Value = 0.2608640200510233
$ws.Range("Q12").Value = 21.00757069772177
$ws.Range("R12").Value = 189.068136279496
$ws.Range("S12").Value = 0.09599936721859989
$ws.Range("T12").Value = 0.0959993672185999

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5451493333333333
$ws.Range("H13").Value = 1.635448
$ws.Range("I13").Value = 0.3680053968340403
$ws.Range("J13").Value = 0.3680053968340404
$ws.Range("M13").Value = 16.07945366666667
$ws.Range("N13").Value = 48.238361
$ws.Range("O13").Value = 0.1088491702632547
$ws.Range("P13").Value = 0.1088491702632547
$ws.Range("Q13").Value = 8.765703446747555
$ws.Range("R13").Value = 78.89133102072799
$ws.Range("S13").Value = 0.04005708209778508
$ws.Range("T13").Value = 0.04005708209778509

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05723766666666667
$ws.Range("H14").Value = 0.171713
$ws.Range("I14").Value = 0.03863853250397663
$ws.Range("J14").Value = 0.03863853250397663
$ws.Range("M14").Value = 46.29121633333333
$ws.Range("N14").Value = 138.873649
$ws.Range("O14").Value = 0.3133663986859022
$ws.Range("P14").Value = 0.3133663986859022
$ws.Range("Q14").Value = 2.649601210081889
$ws.Range("R14").Value = 23.846410890737
$ws.Range("S14").Value = 0.01210801778127933
$ws.Range("T14").Value = 0.01210801778127933

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05723766666666667
$ws.Range("H15").Value = 0.171713
$ws.Range("I15").Value = 0.03863853250397663
$ws.Range("J15").Value = 0.03863853250397663
$ws.Range("M15").Value = 46.81622333333333
$ws.Range("O15").Value = 0.3169204109998198
$ws.Range("P15").Value = 0.3169204109998198
$ws.Range("Q15").Value = 2.679651385745556
$ws.Range("R15").Value = 24.11686247171
$ws.Range("S15").Value = 0.01224533960159017
$ws.Range("T15").Value = 0.01224533960159017

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05723766666666667
$ws.Range("H16").Value = 0.171713
$ws.Range("I16").Value = 0.03863853250397663
$ws.Range("J16").Value = 0.03863853250397663
$ws.Range("M16").Value = 38.53544233333333
$ws.Range("N16").Value = 115.606327
$ws.Range("O16").Value = 0.2608640200510233
$ws.Range("P16").Value = 0.2608640200510233
$ws.Range("Q16").Value = 2.205678803127888
$ws.Range("R16").Value = 19.851109228151
$ws.Range("S16").Value = 0.01007940291785948
$ws.Range("T16").Value = 0.01007940291785948

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05723766666666667
$ws.Range("H17").Value = 0.171713
$ws.Range("I17").Value = 0.03863853250397663
$ws.Range("J17").Value = 0.03863853250397663
$ws.Range("M17").Value = 16.07945366666667
$ws.Range("N17").Value = 48.238361
$ws.Range("O17").Value = 0.1088491702632547
$ws.Range("P17").Value = 0.1088491702632547
$ws.Range("Q17").Value = 0.9203504091547777
$ws.Range("R17").Value = 8.283153682393
$ws.Range("S17").Value = 0.004205772203247654
$ws.Range("T17").Value = 0.004205772203247655
